$d = $word.ActiveDocument

# The last paragraph currently holds the "Sum 1 + 2 + 4 + 8 +... + n "
# run and, at its very end, the _GoBack bookmark. We want to insert four
# new bulleted list items *before* that paragraph's final paragraph mark
# (so the bookmark stays attached to what becomes the new last, empty,
# list paragraph), leaving the "Sum..." paragraph itself untouched.

$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range

# Position right before the paragraph's own end-of-paragraph mark.
$insertPoint = $d.Range($lastRange.End - 1, $lastRange.End - 1)

$newItems = @(
    'Soda can – optimal material usage for the given volume',
    'Some equation has to be solved using iterations (Newton’s tangent method); write in some Python code that completes the task',
    'Write some JavaScript that draws something in 2D. ',
    'Limits, integrals, derivatives, series… Use limit to describe something from the real life (a process that ultimately stabilizes?)'
)

$text = "`r" + ($newItems -join "`r") + "`r"
$insertPoint.InsertBefore($text)

Write-Host "Inserted $($newItems.Count) new list items"
